# "update entregable 1, 2" - refresh the regression-run results on the
# DepositoLocalEfectico sheet with a new test execution's output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 holds the latest test run: Estado (E), Transaccion (F), Fecha (G)
$ws.Range("E2").Value = "PASSED"
$ws.Range("F2").Value = "TT23195K05WV 09:3"
$ws.Range("G2").Value = "14 jul. 2023, 09:39:17"

# Leave the selection where the user's cursor ended up after the edit.
$ws.Range("I5").Select()
